$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the duplicated "Debug" mini-table in columns F, G, H
$ws.Range("F1:H2").ClearContents()

# Rename "Print Plots?" header to "Output Plots?" and update its note
$ws.Range("A10").Value = "Output Plots?"
$ws.Range("D10").Value = "If set to no, nothing in this section matters."

# Update Plot Contours / Plot Hatches values to Yes
$ws.Range("B13").Value = "Yes"
$ws.Range("B14").Value = "Yes"

# Update selection to reflect new active cell
$ws.Range("D11").Select()
